# Update the cryptos price/volume snapshot (GitHub Actions refresh).
# Numeric-looking price strings are prefixed with a leading apostrophe so
# Excel stores them as text (matching the original inline-string cells)
# instead of silently parsing them into numbers/dates.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.413.68"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "1.849.64"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'240.80"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").Value = "'0.6299"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D8").Value = "'0.07688"
$ws.Range("E8").Value = "  +1.90%  "
$ws.Range("D9").Value = "'0.2944"
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("D11").Value = "'0.07749"
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("D12").Value = "1.853.69"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").Value = "'5.023"
$ws.Range("E13").Value = "  +0.57%  "
$ws.Range("D14").Value = "'0.00001083"
$ws.Range("E14").Value = "  +7.98%  "
$ws.Range("D15").Value = "'0.6805"
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("E16").Value = "  +1.02%  "
$ws.Range("D17").Value = "2.103.06"
$ws.Range("E17").Value = "  -0.70%  "
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").Value = "29.426.56"
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("D20").Value = "'229.34"
$ws.Range("E20").Value = "  +0.71%  "
$ws.Range("D21").Value = "'12.48"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "'7.450"
$ws.Range("E23").Value = "  -1.42%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").Value = "'157.26"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  -0.81%  "
$ws.Range("D27").Value = "'8.382"
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("E29").Value = "  +3.94%  "
$ws.Range("D30").Value = "'1.467"
$ws.Range("E30").Value = "  +0.23%  "
$ws.Range("D31").Value = "'0.05718"
$ws.Range("E31").Value = "  +0.67%  "
$ws.Range("D32").Value = "'4.116"
$ws.Range("E32").Value = "  -0.12%  "
$ws.Range("E33").Value = "  +0.88%  "
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("D35").Value = "'1.160"
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("D36").Value = "'0.7094"
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").Value = "'2.777"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D40").Value = "1.224.33"
$ws.Range("D41").Value = "'6.449"
$ws.Range("E41").Value = "  +4.57%  "
$ws.Range("E42").Value = "  +0.21%  "
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").Value = "2.011.91"
$ws.Range("E44").Value = "  -0.72%  "
$ws.Range("D45").Value = "'101.86"
$ws.Range("E45").Value = "  +0.55%  "
$ws.Range("D46").Value = "'66.27"
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("E47").Value = "  +2.69%  "
$ws.Range("E48").Value = "  +0.97%  "
$ws.Range("D49").Value = "'0.4026"
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'9.033"
$ws.Range("E50").Value = "  -0.77%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "'1.687"
$ws.Range("E51").Value = "  +0.18%  "
